$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("C2").Value = 11
$ws.Range("C4").Value = 1.2

# Update the active selection to C4 (as reflected in sheetView)
$ws.Range("C4").Select()
